# Adds two new columns (G = "Org_Size", H = "CC") with their data to the
# existing "Estrutura" table, matching the header/data style of the other
# columns, then re-selects the whole used range (A1:H1048576) as the last
# user action before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header + data values -------------------------------------------------
# Note: "CC" is typed before "Org_Size" so the shared-string table gets the
# two new entries in that order (CC first, Org_Size second).
$ws.Range("H1").Value = "CC"
$ws.Range("G1").Value = "Org_Size"

$ws.Range("G2").Value = 100
$ws.Range("H2").Value = 200

$ws.Range("G3").Value = 1000
$ws.Range("H3").Value = 2000

# --- Match formatting of the existing columns ---------------------------------
# Header style (column F1) -> new header cells G1:H1
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Data style (column F2:F3) -> new data cells G2:H3
$ws.Range("F2:F3").Copy()
$ws.Range("G2:H3").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Size the new columns to fit their content --------------------------------
$ws.Range("H1:H3").ColumnWidth = 4.166666666666667
$ws.Range("G1:G3").ColumnWidth = 7.833333333333333

# --- Re-select the full used range, as the last interactive action -----------
$ws.Range("A1:H1048576").Select()
